# Scheduled-runner update: recalculated profit figures across all profession sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per "chore: update Sheets via scheduled runner".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 159701.84
$ws.Range("I17").Value = 1400
$ws.Range("J17").Value = 163980.27
$ws.Range("K17").Value = 4200
$ws.Range("L17").Value = 491940.8099999999
$ws.Range("M17").Value = -4032
$ws.Range("N17").Value = -492276.8099999999
$ws.Range("H33").Value = 27778208
$ws.Range("I33").Value = 228.96153
$ws.Range("J33").Value = 100000950
$ws.Range("K33").Value = 228.96153
$ws.Range("L33").Value = 100000950
$ws.Range("M33").Value = 0.03846999999998957
$ws.Range("N33").Value = -100001408
$ws.Range("H111").Value = 1071.4286
$ws.Range("I111").Value = 666.6667
$ws.Range("J111").Value = 1375
$ws.Range("K111").Value = 2000.0001
$ws.Range("L111").Value = 4125
$ws.Range("M111").Value = 1066.9999
$ws.Range("N111").Value = -10259
$ws.Range("H116").Value = 2974.75
$ws.Range("J116").Value = 2974.75
$ws.Range("L116").Value = 2974.75
$ws.Range("N116").Value = -9858.75
$ws.Range("H132").Value = 4669.839
$ws.Range("I132").Value = 1437.75
$ws.Range("J132").Value = 15751.286
$ws.Range("K132").Value = 4313.25
$ws.Range("L132").Value = 47253.858
$ws.Range("M132").Value = -1783.25
$ws.Range("N132").Value = -52313.858
$ws.Range("H137").Value = 1815.2646
$ws.Range("I137").Value = 1406.6666
$ws.Range("J137").Value = 2274.9375
$ws.Range("K137").Value = 4219.9998
$ws.Range("L137").Value = 6824.8125
$ws.Range("M137").Value = -1669.9998
$ws.Range("N137").Value = -11924.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3673.9
$ws.Range("I32").Value = 3258.152
$ws.Range("J32").Value = 8455
$ws.Range("K32").Value = 3258.152
$ws.Range("L32").Value = 8455
$ws.Range("M32").Value = -2971.152
$ws.Range("N32").Value = -9029
$ws.Range("H74").Value = 91784
$ws.Range("I74").Value = 91784
$ws.Range("K74").Value = 91784
$ws.Range("M74").Value = -90910
$ws.Range("H77").Value = 91784
$ws.Range("I77").Value = 91784
$ws.Range("K77").Value = 458920
$ws.Range("M77").Value = -454552
$ws.Range("H88").Value = 1784.6471
$ws.Range("I88").Value = 1591.4286
$ws.Range("J88").Value = 1919.9
$ws.Range("K88").Value = 1591.4286
$ws.Range("L88").Value = 1919.9
$ws.Range("M88").Value = -1185.4286
$ws.Range("N88").Value = -2731.9
$ws.Range("H91").Value = 1784.6471
$ws.Range("I91").Value = 1591.4286
$ws.Range("J91").Value = 1919.9
$ws.Range("K91").Value = 1591.4286
$ws.Range("L91").Value = 1919.9
$ws.Range("M91").Value = -187.4286
$ws.Range("N91").Value = -4727.9
$ws.Range("H132").Value = 1610
$ws.Range("I132").Value = 1176.5483
$ws.Range("J132").Value = 3103
$ws.Range("K132").Value = 3529.6449
$ws.Range("L132").Value = 9309
$ws.Range("M132").Value = -999.6448999999998
$ws.Range("N132").Value = -14369

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 7641.4287
$ws.Range("I5").Value = 3833.3333
$ws.Range("J5").Value = 10497.5
$ws.Range("K5").Value = 3833.3333
$ws.Range("L5").Value = 10497.5
$ws.Range("M5").Value = -3720.3333
$ws.Range("N5").Value = -10723.5
$ws.Range("H126").Value = 33772.188
$ws.Range("J126").Value = 33772.188
$ws.Range("L126").Value = 33772.188
$ws.Range("N126").Value = -43652.188
$ws.Range("H134").Value = 669851.2
$ws.Range("I134").Value = 836095.7
$ws.Range("J134").Value = 4873
$ws.Range("K134").Value = 2508287.1
$ws.Range("L134").Value = 14619
$ws.Range("M134").Value = -2505752.1
$ws.Range("N134").Value = -19689

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4861.2188
$ws.Range("I58").Value = 8208.467000000001
$ws.Range("J58").Value = 1907.7646
$ws.Range("K58").Value = 8208.467000000001
$ws.Range("L58").Value = 1907.7646
$ws.Range("M58").Value = -8005.467000000001
$ws.Range("N58").Value = -2313.7646
$ws.Range("H74").Value = 17578.5
$ws.Range("J74").Value = 17578.5
$ws.Range("L74").Value = 17578.5
$ws.Range("N74").Value = -19326.5
$ws.Range("H77").Value = 17578.5
$ws.Range("J77").Value = 17578.5
$ws.Range("L77").Value = 52735.5
$ws.Range("N77").Value = -61471.5
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 1123785.8
$ws.Range("I132").Value = 1259.84
$ws.Range("J132").Value = 4631679
$ws.Range("K132").Value = 3779.52
$ws.Range("L132").Value = 13895037
$ws.Range("M132").Value = -1249.52
$ws.Range("N132").Value = -13900097
$ws.Range("H134").Value = 2741.7144
$ws.Range("I134").Value = 1667.7693
$ws.Range("K134").Value = 5003.3079
$ws.Range("M134").Value = -2468.3079
$ws.Range("H136").Value = 4861.2188
$ws.Range("I136").Value = 8208.467000000001
$ws.Range("J136").Value = 1907.7646
$ws.Range("K136").Value = 24625.401
$ws.Range("L136").Value = 5723.293799999999
$ws.Range("M136").Value = -22075.401
$ws.Range("N136").Value = -10823.2938
$ws.Range("H141").Value = 56081
$ws.Range("J141").Value = 56081
$ws.Range("L141").Value = 56081
$ws.Range("N141").Value = -66441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1276.2693
$ws.Range("I131").Value = 5450
$ws.Range("J131").Value = 1166.4342
$ws.Range("K131").Value = 16350
$ws.Range("L131").Value = 3499.3026
$ws.Range("M131").Value = -11310
$ws.Range("N131").Value = -13579.3026

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 798.3
$ws.Range("I107").Value = 716.6667
$ws.Range("J107").Value = 920.75
$ws.Range("K107").Value = 716.6667
$ws.Range("L107").Value = 920.75
$ws.Range("M107").Value = 1203.3333
$ws.Range("N107").Value = -4760.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5023.148
$ws.Range("I136").Value = 1362.7778
$ws.Range("J136").Value = 12343.889
$ws.Range("K136").Value = 4088.3334
$ws.Range("L136").Value = 37031.667
$ws.Range("M136").Value = -1538.3334
$ws.Range("N136").Value = -42131.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 23998
$ws.Range("I106").Value = 23998
$ws.Range("K106").Value = 23998
$ws.Range("M106").Value = -22736
$ws.Range("H132").Value = 3259.5193
$ws.Range("I132").Value = 3473.1875
$ws.Range("J132").Value = 2917.65
$ws.Range("K132").Value = 10419.5625
$ws.Range("L132").Value = 8752.950000000001
$ws.Range("M132").Value = -7889.5625
$ws.Range("N132").Value = -13812.95
$ws.Range("H136").Value = 4287.34
$ws.Range("I136").Value = 2283.394
$ws.Range("J136").Value = 8177.353
$ws.Range("K136").Value = 6850.181999999999
$ws.Range("L136").Value = 24532.059
$ws.Range("M136").Value = -4300.181999999999
$ws.Range("N136").Value = -29632.059
